$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.295571327209473
$ws.Range("B1").Value = 3.839211940765381
$ws.Range("C1").Value = 3.936913013458252
$ws.Range("D1").Value = 1.726462721824646
$ws.Range("E1").Value = 1.213523626327515
